$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.372558333333333
$ws.Range("H2").Value = 4.117675
$ws.Range("I2").Value = 0.3340102211301095
$ws.Range("J2").Value = 0.3340102211301095
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 2.896460896252778
$ws.Range("R2").Value = 26.068148066275
$ws.Range("S2").Value = 0.1279940923657414
$ws.Range("T2").Value = 0.1279940923657414

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.372558333333333
$ws.Range("H3").Value = 4.117675
$ws.Range("I3").Value = 0.3340102211301095
$ws.Range("J3").Value = 0.3340102211301095
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("Q3").Value = 0.7151760947777777
$ws.Range("R3").Value = 6.436584852999999
$ws.Range("S3").Value = 0.03160350455660647
$ws.Range("T3").Value = 0.03160350455660647

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.372558333333333
$ws.Range("H4").Value = 4.117675
$ws.Range("I4").Value = 0.3340102211301095
$ws.Range("J4").Value = 0.3340102211301095
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 2.834865968486111
$ws.Range("R4").Value = 25.513793716375
$ws.Range("S4").Value = 0.1252722234518447
$ws.Range("T4").Value = 0.1252722234518447

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.372558333333333
$ws.Range("H5").Value = 4.117675
$ws.Range("I5").Value = 0.3340102211301095
$ws.Range("J5").Value = 0.3340102211301095
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 1.112029833447222
$ws.Range("R5").Value = 10.008268501025
$ws.Range("S5").Value = 0.04914040075591693
$ws.Range("T5").Value = 0.04914040075591693

# Row 6
$ws.Range("I6").Value = 0.01293592767872722
$ws.Range("J6").Value = 0.01293592767872721
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.110264333333333
$ws.Range("N6").Value = 6.330793
$ws.Range("O6").Value = 0.3832041185227171
$ws.Range("P6").Value = 0.3832041185227171
$ws.Range("Q6").Value = 0.1121774314313333
$ws.Range("R6").Value = 1.009596882882
$ws.Range("S6").Value = 0.00495710076340028
$ws.Range("T6").Value = 0.00495710076340028

# Row 7
$ws.Range("I7").Value = 0.01293592767872722
$ws.Range("J7").Value = 0.01293592767872721
$ws.Range("O7").Value = 0.0946183755984393
$ws.Range("P7").Value = 0.0946183755984393
$ws.Range("S7").Value = 0.001223976463820059
$ws.Range("T7").Value = 0.001223976463820059

# Row 8
$ws.Range("I8").Value = 0.01293592767872722
$ws.Range("J8").Value = 0.01293592767872721
$ws.Range("M8").Value = 2.065388333333333
$ws.Range("N8").Value = 6.196165
$ws.Range("O8").Value = 0.3750550597762889
$ws.Range("P8").Value = 0.3750550597762889
$ws.Range("Q8").Value = 0.1097919130233333
$ws.Range("R8").Value = 0.98812721721
$ws.Range("S8").Value = 0.004851685128806786
$ws.Range("T8").Value = 0.004851685128806785

# Row 9
$ws.Range("I9").Value = 0.01293592767872722
$ws.Range("J9").Value = 0.01293592767872721
$ws.Range("M9").Value = 0.8101876666666666
$ws.Range("N9").Value = 2.430563
$ws.Range("O9").Value = 0.1471224461025547
$ws.Range("P9").Value = 0.1471224461025547
$ws.Range("Q9").Value = 0.04306795598466667
$ws.Range("R9").Value = 0.387611603862
$ws.Range("S9").Value = 0.001903165322700091
$ws.Range("T9").Value = 0.001903165322700091

# Row 10
$ws.Range("G10").Value = 2.683614
$ws.Range("H10").Value = 8.050841999999999
$ws.Range("I10").Value = 0.6530538511911632
$ws.Range("J10").Value = 0.6530538511911632
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.110264333333333
$ws.Range("N10").Value = 6.330793
$ws.Range("O10").Value = 0.3832041185227171
$ws.Range("P10").Value = 0.3832041185227171
$ws.Range("Q10").Value = 5.663134908633999
$ws.Range("R10").Value = 50.968214177706
$ws.Range("S10").Value = 0.2502529253935754
$ws.Range("T10").Value = 0.2502529253935754

# Row 11
$ws.Range("G11").Value = 2.683614
$ws.Range("H11").Value = 8.050841999999999
$ws.Range("I11").Value = 0.6530538511911632
$ws.Range("J11").Value = 0.6530538511911632
$ws.Range("O11").Value = 0.0946183755984393
$ws.Range("P11").Value = 0.0946183755984393
$ws.Range("Q11").Value = 1.39830602008
$ws.Range("R11").Value = 12.58475418072
$ws.Range("S11").Value = 0.06179089457801277
$ws.Range("T11").Value = 0.06179089457801277

# Row 12
$ws.Range("G12").Value = 2.683614
$ws.Range("H12").Value = 8.050841999999999
$ws.Range("I12").Value = 0.6530538511911632
$ws.Range("J12").Value = 0.6530538511911632
$ws.Range("M12").Value = 2.065388333333333
$ws.Range("N12").Value = 6.196165
$ws.Range("O12").Value = 0.3750550597762889
$ws.Range("P12").Value = 0.3750550597762889
$ws.Range("Q12").Value = 5.54270504677
$ws.Range("R12").Value = 49.88434542093
$ws.Range("S12").Value = 0.2449311511956374
$ws.Range("T12").Value = 0.2449311511956374

# Row 13
$ws.Range("G13").Value = 2.683614
$ws.Range("H13").Value = 8.050841999999999
$ws.Range("I13").Value = 0.6530538511911632
$ws.Range("J13").Value = 0.6530538511911632
$ws.Range("M13").Value = 0.8101876666666666
$ws.Range("N13").Value = 2.430563
$ws.Range("O13").Value = 0.1471224461025547
$ws.Range("P13").Value = 0.1471224461025547
$ws.Range("Q13").Value = 2.174230964894
$ws.Range("R13").Value = 19.568078684046
$ws.Range("S13").Value = 0.0960788800239377
$ws.Range("T13").Value = 0.0960788800239377
